# Insert a new "Visible" column (G) on the Input sheet, between "Enabled"
# and "Options", filled with TRUE for every data row, and make the Input
# sheet the active sheet/selection (instead of Output).

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")

# Shift the old "Options"/"Errormessage" columns one to the right and
# insert the new "Visible" column in their place.
$wsInput.Range("G1").EntireColumn.Insert()

$wsInput.Range("G1").Value = "Visible"
$wsInput.Range("G2").Value = $true
$wsInput.Range("G3").Value = $true
$wsInput.Range("G4").Value = $true
$wsInput.Range("G5").Value = $true

# Input becomes the active sheet/tab, with G6 selected.
$wsInput.Activate()
$wsInput.Range("G6").Select()
